$d = $word.ActiveDocument
$full = $d.Content.Text
$idx = $full.IndexOf("Sukabumi,")
$len = "Sukabumi,".Length
$r1 = $d.Range($idx, $idx+$len)
$r1.Text = "{"
$full2 = $d.Content.Text
$idx2 = $full2.IndexOf("{{tanggal_surat}")
$r2 = $d.Range($idx2+1, $idx2+1)
try {
  $d.Bookmarks.Add("split_marker", $r2)
  Write-Output "bm added"
} catch {
  Write-Output "ERR: $_"
}
$r3 = $d.Range($idx2, $idx2+1)
$r3.Delete()
try {
  $d.Bookmarks("split_marker").Delete()
  Write-Output "bm deleted"
} catch {
  Write-Output "ERR2: $_"
}
Write-Output $d.Content.Text
